$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.16
$ws.Range("C2").Value = 0.62
$ws.Range("J2").Value = 0.01
$ws.Range("P2").Value = 0.1066666666666667
$ws.Range("S2").Value = 0.1033333333333333
$ws.Range("B3").Value = 0.0101010101010101
$ws.Range("C3").Value = 0.02525252525252525
$ws.Range("J3").Value = 0.04545454545454546
$ws.Range("P3").Value = 0.7424242424242424
$ws.Range("S3").Value = 0.1767676767676768
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.6470588235294118
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.08333333333333333
$ws.Range("D6").Value = 0.007936507936507936
$ws.Range("F6").Value = 0.09126984126984126
$ws.Range("J6").Value = 0.2063492063492063
$ws.Range("O6").Value = 0.01984126984126984
$ws.Range("Q6").Value = 0.123015873015873
$ws.Range("R6").Value = 0.09523809523809523
$ws.Range("S6").Value = 0.373015873015873
$ws.Range("B7").Value = 0.1216216216216216
$ws.Range("D7").Value = 0.02027027027027027
$ws.Range("F7").Value = 0.06081081081081081
$ws.Range("J7").Value = 0.1013513513513514
$ws.Range("O7").Value = 0.02027027027027027
$ws.Range("Q7").Value = 0.1621621621621622
$ws.Range("R7").Value = 0.1081081081081081
$ws.Range("S7").Value = 0.4054054054054054
$ws.Range("B8").Value = 0.1129807692307692
$ws.Range("D8").Value = 0.02403846153846154
$ws.Range("E8").Value = 0.002403846153846154
$ws.Range("F8").Value = 0.08173076923076923
$ws.Range("J8").Value = 0.1153846153846154
$ws.Range("O8").Value = 0.007211538461538462
$ws.Range("Q8").Value = 0.1730769230769231
$ws.Range("R8").Value = 0.06971153846153846
$ws.Range("S8").Value = 0.4134615384615384
$ws.Range("B9").Value = 0.12
$ws.Range("D9").Value = 0.02
$ws.Range("F9").Value = 0.125
$ws.Range("J9").Value = 0.095
$ws.Range("O9").Value = 0.01
$ws.Range("Q9").Value = 0.19
$ws.Range("R9").Value = 0.07000000000000001
$ws.Range("S9").Value = 0.37
$ws.Range("B10").Value = 0.1153519932145886
$ws.Range("D10").Value = 0.01611535199321459
$ws.Range("E10").Value = 0.0008481764206955047
$ws.Range("F10").Value = 0.07718405428329092
$ws.Range("J10").Value = 0.1128074639525021
$ws.Range("O10").Value = 0.01357082273112807
$ws.Range("Q10").Value = 0.1959287531806616
$ws.Range("R10").Value = 0.08566581849024597
$ws.Range("S10").Value = 0.3825275657336726
$ws.Range("G11").Value = 0.1403508771929824
$ws.Range("J11").Value = 0.08771929824561403
$ws.Range("K11").Value = 0.2105263157894737
$ws.Range("L11").Value = 0.5228070175438596
$ws.Range("S11").Value = 0.03859649122807018
$ws.Range("G12").Value = 0.5949367088607594
$ws.Range("J12").Value = 0.2721518987341772
$ws.Range("K12").Value = 0.02531645569620253
$ws.Range("L12").Value = 0.03164556962025317
$ws.Range("S12").Value = 0.0759493670886076
$ws.Range("F15").Value = 0.02116402116402116
$ws.Range("H15").Value = 0.1164021164021164
$ws.Range("I15").Value = 0.06349206349206349
$ws.Range("J15").Value = 0.4656084656084656
$ws.Range("K15").Value = 0.08994708994708994
$ws.Range("M15").Value = 0.01587301587301587
$ws.Range("O15").Value = 0.04232804232804233
$ws.Range("S15").Value = 0.1851851851851852
$ws.Range("F16").Value = 0.005235602094240838
$ws.Range("H16").Value = 0.1465968586387434
$ws.Range("I16").Value = 0.07329842931937172
$ws.Range("J16").Value = 0.4083769633507853
$ws.Range("K16").Value = 0.1465968586387434
$ws.Range("M16").Value = 0.02617801047120419
$ws.Range("O16").Value = 0.05759162303664921
$ws.Range("S16").Value = 0.1361256544502618
$ws.Range("F17").Value = 0.02040816326530612
$ws.Range("H17").Value = 0.2244897959183673
$ws.Range("I17").Value = 0.1147959183673469
$ws.Range("J17").Value = 0.4285714285714285
$ws.Range("K17").Value = 0.07142857142857142
$ws.Range("M17").Value = 0.00510204081632653
$ws.Range("O17").Value = 0.04336734693877551
$ws.Range("S17").Value = 0.09183673469387756
$ws.Range("F18").Value = 0.02173913043478261
$ws.Range("H18").Value = 0.1684782608695652
$ws.Range("I18").Value = 0.09782608695652174
$ws.Range("J18").Value = 0.358695652173913
$ws.Range("K18").Value = 0.1304347826086956
$ws.Range("M18").Value = 0.005434782608695652
$ws.Range("N18").Value = 0.005434782608695652
$ws.Range("O18").Value = 0.07065217391304347
$ws.Range("S18").Value = 0.1413043478260869
$ws.Range("F19").Value = 0.01957585644371941
$ws.Range("H19").Value = 0.2030995106035889
$ws.Range("I19").Value = 0.09053833605220228
$ws.Range("J19").Value = 0.3588907014681892
$ws.Range("K19").Value = 0.1019575856443719
$ws.Range("M19").Value = 0.02039151712887439
$ws.Range("N19").Value = 0.0008156606851549756
$ws.Range("O19").Value = 0.07340946166394779
$ws.Range("S19").Value = 0.1313213703099511
